$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2128.3333
$ws.Range("I18").Value = 1842.5
$ws.Range("J18").Value = 2700
$ws.Range("K18").Value = 1842.5
$ws.Range("L18").Value = 2700
$ws.Range("M18").Value = -1558.5
$ws.Range("N18").Value = -3268

$ws.Range("H74").Value = 3975.2964
$ws.Range("I74").Value = 3932.5557
$ws.Range("J74").Value = 3996.6667
$ws.Range("K74").Value = 3932.5557
$ws.Range("L74").Value = 3996.6667
$ws.Range("M74").Value = -2996.5557
$ws.Range("N74").Value = -5868.6667

$ws.Range("H77").Value = 3975.2964
$ws.Range("I77").Value = 3932.5557
$ws.Range("J77").Value = 3996.6667
$ws.Range("K77").Value = 19662.7785
$ws.Range("L77").Value = 19983.3335
$ws.Range("M77").Value = -14982.7785
$ws.Range("N77").Value = -29343.3335

$ws.Range("H125").Value = 1068.5
$ws.Range("I125").Value = 728
$ws.Range("K125").Value = 6552
$ws.Range("M125").Value = -4092

$ws.Range("H137").Value = 2367.3572
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 2367.3572
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 7102.071599999999
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -12202.0716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20943.146
$ws.Range("I32").Value = 18482.967
$ws.Range("J32").Value = 39394.5
$ws.Range("K32").Value = 18482.967
$ws.Range("L32").Value = 39394.5
$ws.Range("M32").Value = -18195.967
$ws.Range("N32").Value = -39968.5

$ws.Range("H61").Value = 105369750
$ws.Range("I61").Value = 66734696
$ws.Range("K61").Value = 66734696
$ws.Range("M61").Value = -66734484

$ws.Range("H136").Value = 105369750
$ws.Range("I136").Value = 66734696
$ws.Range("K136").Value = 200204088
$ws.Range("M136").Value = -200201538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1660.25
$ws.Range("I94").Value = 703.5
$ws.Range("J94").Value = 2617
$ws.Range("K94").Value = 703.5
$ws.Range("L94").Value = 2617
$ws.Range("M94").Value = -252.5
$ws.Range("N94").Value = -3519

$ws.Range("H134").Value = 2209.6
$ws.Range("I134").Value = 2269.92
$ws.Range("J134").Value = 1908
$ws.Range("K134").Value = 6809.76
$ws.Range("L134").Value = 5724
$ws.Range("M134").Value = -4274.76
$ws.Range("N134").Value = -10794

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1378
$ws.Range("I10").Value = 1378
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1378
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -1239
$ws.Range("N10").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 26400.375
$ws.Range("I137").Value = 901.5833
$ws.Range("J137").Value = 51899.168
$ws.Range("K137").Value = 2704.7499
$ws.Range("L137").Value = 155697.504
$ws.Range("M137").Value = 2395.2501
$ws.Range("N137").Value = -165897.504

$ws.Range("H140").Value = 2094.8462
$ws.Range("I140").Value = 1623.4783
$ws.Range("J140").Value = 2291.9636
$ws.Range("K140").Value = 4870.4349
$ws.Range("L140").Value = 6875.8908
$ws.Range("M140").Value = 309.5650999999998
$ws.Range("N140").Value = -17235.8908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1590.6471
$ws.Range("I97").Value = 1534.6666
$ws.Range("J97").Value = 2010.5
$ws.Range("K97").Value = 1534.6666
$ws.Range("L97").Value = 2010.5
$ws.Range("M97").Value = -1038.6666
$ws.Range("N97").Value = -3002.5

$ws.Range("H102").Value = 2400.8572
$ws.Range("I102").Value = 2255.6365
$ws.Range("J102").Value = 2933.3333
$ws.Range("K102").Value = 2255.6365
$ws.Range("L102").Value = 2933.3333
$ws.Range("M102").Value = -633.6365000000001
$ws.Range("N102").Value = -6177.3333

$ws.Range("H132").Value = 97169.95
$ws.Range("I132").Value = 143810.28
$ws.Range("J132").Value = 73849.78999999999
$ws.Range("K132").Value = 431430.84
$ws.Range("L132").Value = 221549.37
$ws.Range("M132").Value = -428900.84
$ws.Range("N132").Value = -226609.37

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2075.3333
$ws.Range("I16").Value = 1045.8823
$ws.Range("J16").Value = 6450.5
$ws.Range("K16").Value = 1045.8823
$ws.Range("L16").Value = 6450.5
$ws.Range("M16").Value = -875.8823
$ws.Range("N16").Value = -6790.5

$ws.Range("H46").Value = 1329.1428
$ws.Range("I46").Value = 1100
$ws.Range("J46").Value = 1420.8
$ws.Range("K46").Value = 1100
$ws.Range("L46").Value = 1420.8
$ws.Range("M46").Value = -912
$ws.Range("N46").Value = -1796.8

$ws.Range("H61").Value = 1909.421
$ws.Range("I61").Value = 1922.3529
$ws.Range("J61").Value = 1799.5
$ws.Range("K61").Value = 1922.3529
$ws.Range("L61").Value = 1799.5
$ws.Range("M61").Value = -1720.3529
$ws.Range("N61").Value = -2203.5

$ws.Range("H93").Value = 3000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 3000
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -5496

$ws.Range("H98").Value = 29000
$ws.Range("J98").Value = 29000
$ws.Range("L98").Value = 29000
$ws.Range("N98").Value = -34990

$ws.Range("H113").Value = 1909.421
$ws.Range("I113").Value = 1922.3529
$ws.Range("J113").Value = 1799.5
$ws.Range("K113").Value = 1922.3529
$ws.Range("L113").Value = 1799.5
$ws.Range("M113").Value = 247.6470999999999
$ws.Range("N113").Value = -6139.5

$ws.Range("H132").Value = 48403.09
$ws.Range("I132").Value = 1621.5555
$ws.Range("J132").Value = 80790.30499999999
$ws.Range("K132").Value = 4864.666499999999
$ws.Range("L132").Value = 242370.915
$ws.Range("M132").Value = -2334.666499999999
$ws.Range("N132").Value = -247430.915

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 31266.666
$ws.Range("J112").Value = 31266.666
$ws.Range("L112").Value = 31266.666
$ws.Range("N112").Value = -34220.666

$ws.Range("H126").Value = 1008.2353
$ws.Range("I126").Value = 883.0968
$ws.Range("J126").Value = 2301.3333
$ws.Range("K126").Value = 2649.2904
$ws.Range("L126").Value = 6903.999899999999
$ws.Range("M126").Value = -179.2903999999999
$ws.Range("N126").Value = -11843.9999

$ws.Range("H132").Value = 72921.57000000001
$ws.Range("I132").Value = 53639.26
$ws.Range("J132").Value = 113628.664
$ws.Range("K132").Value = 160917.78
$ws.Range("L132").Value = 340885.992
$ws.Range("M132").Value = -158387.78
$ws.Range("N132").Value = -345945.992

$ws.Range("H135").Value = 43249.5
$ws.Range("J135").Value = 43249.5
$ws.Range("L135").Value = 43249.5
$ws.Range("N135").Value = -53389.5

$ws.Range("H136").Value = 40850.41
$ws.Range("I136").Value = 23253.021
$ws.Range("K136").Value = 69759.06299999999
$ws.Range("M136").Value = -67209.06299999999
